$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LG생활건강")
$ws.Cells.Item(210,4).Value = 780540
$ws.Cells.Item(210,5).Value = 15302459
$ws.Cells.Item(210,6).Value = 223234440000
$ws.Cells.Item(210,7).Value = 4376503274000
$ws.Cells.Item(210,8).Value = 5.099999904632568
$ws.Cells.Item(211,4).Value = 784682
$ws.Cells.Item(211,5).Value = 15302459
$ws.Cells.Item(211,6).Value = 224419052000
$ws.Cells.Item(211,7).Value = 4376503274000
$ws.Cells.Item(211,8).Value = 5.130000114440918
$ws.Cells.Item(212,1).Value = 45968
$ws.Cells.Item(212,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212,2).Value = 289000
$ws.Cells.Item(212,3).Value = 78161
$ws.Cells.Item(212,4).Value = 778365
$ws.Cells.Item(212,5).Value = 15302459
$ws.Cells.Item(212,6).Value = 224947485000
$ws.Cells.Item(212,7).Value = 4422410651000
$ws.Cells.Item(212,8).Value = 5.090000152587891
$ws.Cells.Item(213,1).Value = 45971
$ws.Cells.Item(213,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,2).Value = 288500
$ws.Cells.Item(213,3).Value = 44780
$ws.Cells.Item(213,4).Value = 786346
$ws.Cells.Item(213,5).Value = 15302459
$ws.Cells.Item(213,6).Value = 226860821000
$ws.Cells.Item(213,7).Value = 4414759421500
$ws.Cells.Item(213,8).Value = 5.139999866485596
$ws.Cells.Item(214,1).Value = 45972
$ws.Cells.Item(214,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214,2).Value = 291000
$ws.Cells.Item(214,3).Value = 83795
$ws.Cells.Item(214,4).Value = 769888
$ws.Cells.Item(214,5).Value = 15302459
$ws.Cells.Item(214,6).Value = 224037408000
$ws.Cells.Item(214,7).Value = 4453015569000
$ws.Cells.Item(214,8).Value = 5.03000020980835
$ws.Cells.Item(215,1).Value = 45973
$ws.Cells.Item(215,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215,2).Value = 298500
$ws.Cells.Item(215,3).Value = 121564
$ws.Cells.Item(215,4).Value = 709982
$ws.Cells.Item(215,5).Value = 15302459
$ws.Cells.Item(215,6).Value = 211929627000
$ws.Cells.Item(215,7).Value = 4567784011500
$ws.Cells.Item(215,8).Value = 4.639999866485596
$ws.Cells.Item(216,1).Value = 45974
$ws.Cells.Item(216,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216,2).Value = 301000
$ws.Cells.Item(216,3).Value = 76998
$ws.Cells.Item(217,1).Value = 45975
$ws.Cells.Item(217,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217,2).Value = 294000
$ws.Cells.Item(217,3).Value = 53460

$ws = $wb.Worksheets.Item("아모레퍼시픽")
$ws.Cells.Item(210,4).Value = 792762
$ws.Cells.Item(210,5).Value = 58492759
$ws.Cells.Item(210,6).Value = 95210716200
$ws.Cells.Item(210,7).Value = 7024980355900
$ws.Cells.Item(210,8).Value = 1.360000014305115
$ws.Cells.Item(211,4).Value = 824440
$ws.Cells.Item(211,5).Value = 58492759
$ws.Cells.Item(211,6).Value = 97778584000
$ws.Cells.Item(211,7).Value = 6937241217400
$ws.Cells.Item(211,8).Value = 1.409999966621399
$ws.Cells.Item(212,1).Value = 45968
$ws.Cells.Item(212,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212,2).Value = 126500
$ws.Cells.Item(212,3).Value = 1604910
$ws.Cells.Item(212,4).Value = 599752
$ws.Cells.Item(212,5).Value = 58492759
$ws.Cells.Item(212,6).Value = 75868628000
$ws.Cells.Item(212,7).Value = 7399334013500
$ws.Cells.Item(212,8).Value = 1.029999971389771
$ws.Cells.Item(213,1).Value = 45971
$ws.Cells.Item(213,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,2).Value = 126100
$ws.Cells.Item(213,3).Value = 299775
$ws.Cells.Item(213,4).Value = 602095
$ws.Cells.Item(213,5).Value = 58492759
$ws.Cells.Item(213,6).Value = 75924179500
$ws.Cells.Item(213,7).Value = 7375936909900
$ws.Cells.Item(213,8).Value = 1.029999971389771
$ws.Cells.Item(214,1).Value = 45972
$ws.Cells.Item(214,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214,2).Value = 123900
$ws.Cells.Item(214,3).Value = 291886
$ws.Cells.Item(214,4).Value = 602865
$ws.Cells.Item(214,5).Value = 58492759
$ws.Cells.Item(214,6).Value = 74694973500
$ws.Cells.Item(214,7).Value = 7247252840100
$ws.Cells.Item(214,8).Value = 1.029999971389771
$ws.Cells.Item(215,1).Value = 45973
$ws.Cells.Item(215,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215,2).Value = 126300
$ws.Cells.Item(215,3).Value = 199969
$ws.Cells.Item(215,4).Value = 556494
$ws.Cells.Item(215,5).Value = 58492759
$ws.Cells.Item(215,6).Value = 70285192200
$ws.Cells.Item(215,7).Value = 7387635461700
$ws.Cells.Item(215,8).Value = 0.949999988079071
$ws.Cells.Item(216,1).Value = 45974
$ws.Cells.Item(216,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216,2).Value = 127800
$ws.Cells.Item(216,3).Value = 271643
$ws.Cells.Item(217,1).Value = 45975
$ws.Cells.Item(217,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217,2).Value = 127400
$ws.Cells.Item(217,3).Value = 333741

$ws = $wb.Worksheets.Item("한국콜마")
$ws.Cells.Item(210,4).Value = 602780
$ws.Cells.Item(210,5).Value = 23605077
$ws.Cells.Item(210,6).Value = 44907110000
$ws.Cells.Item(210,7).Value = 1758578236500
$ws.Cells.Item(210,8).Value = 2.549999952316284
$ws.Cells.Item(211,4).Value = 617784
$ws.Cells.Item(211,5).Value = 23605077
$ws.Cells.Item(211,6).Value = 44974675200
$ws.Cells.Item(211,7).Value = 1718449605600
$ws.Cells.Item(211,8).Value = 2.619999885559082
$ws.Cells.Item(212,1).Value = 45968
$ws.Cells.Item(212,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212,2).Value = 72100
$ws.Cells.Item(212,3).Value = 390579
$ws.Cells.Item(212,4).Value = 628127
$ws.Cells.Item(212,5).Value = 23605077
$ws.Cells.Item(212,6).Value = 45287956700
$ws.Cells.Item(212,7).Value = 1701926051700
$ws.Cells.Item(212,8).Value = 2.660000085830688
$ws.Cells.Item(213,1).Value = 45971
$ws.Cells.Item(213,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,2).Value = 66000
$ws.Cells.Item(213,3).Value = 645826
$ws.Cells.Item(213,4).Value = 670939
$ws.Cells.Item(213,5).Value = 23605077
$ws.Cells.Item(213,6).Value = 44281974000
$ws.Cells.Item(213,7).Value = 1557935082000
$ws.Cells.Item(213,8).Value = 2.839999914169312
$ws.Cells.Item(214,1).Value = 45972
$ws.Cells.Item(214,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214,2).Value = 63600
$ws.Cells.Item(214,3).Value = 699012
$ws.Cells.Item(214,4).Value = 634230
$ws.Cells.Item(214,5).Value = 23605077
$ws.Cells.Item(214,6).Value = 40337028000
$ws.Cells.Item(214,7).Value = 1501282897200
$ws.Cells.Item(214,8).Value = 2.690000057220459
$ws.Cells.Item(215,1).Value = 45973
$ws.Cells.Item(215,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215,2).Value = 64800
$ws.Cells.Item(215,3).Value = 244709
$ws.Cells.Item(215,4).Value = 581644
$ws.Cells.Item(215,5).Value = 23605077
$ws.Cells.Item(215,6).Value = 37690531200
$ws.Cells.Item(215,7).Value = 1529608989600
$ws.Cells.Item(215,8).Value = 2.460000038146973
$ws.Cells.Item(216,1).Value = 45974
$ws.Cells.Item(216,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216,2).Value = 65400
$ws.Cells.Item(216,3).Value = 213672
$ws.Cells.Item(217,1).Value = 45975
$ws.Cells.Item(217,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217,2).Value = 64200
$ws.Cells.Item(217,3).Value = 185901

$ws = $wb.Worksheets.Item("코스맥스")
$ws.Cells.Item(210,4).Value = 471295
$ws.Cells.Item(210,5).Value = 11349509
$ws.Cells.Item(210,6).Value = 92185302000
$ws.Cells.Item(210,7).Value = 2219963960400
$ws.Cells.Item(210,8).Value = 4.150000095367432
$ws.Cells.Item(211,4).Value = 464293
$ws.Cells.Item(211,5).Value = 11349509
$ws.Cells.Item(211,6).Value = 88958538800
$ws.Cells.Item(211,7).Value = 2174565924400
$ws.Cells.Item(211,8).Value = 4.090000152587891
$ws.Cells.Item(212,1).Value = 45968
$ws.Cells.Item(212,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212,2).Value = 193600
$ws.Cells.Item(212,3).Value = 106665
$ws.Cells.Item(212,4).Value = 462283
$ws.Cells.Item(212,5).Value = 11349509
$ws.Cells.Item(212,6).Value = 89497988800
$ws.Cells.Item(212,7).Value = 2197264942400
$ws.Cells.Item(212,8).Value = 4.070000171661377
$ws.Cells.Item(213,1).Value = 45971
$ws.Cells.Item(213,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,2).Value = 189500
$ws.Cells.Item(213,3).Value = 120630
$ws.Cells.Item(213,4).Value = 455981
$ws.Cells.Item(213,5).Value = 11349509
$ws.Cells.Item(213,6).Value = 86408399500
$ws.Cells.Item(213,7).Value = 2150731955500
$ws.Cells.Item(213,8).Value = 4.019999980926514
$ws.Cells.Item(214,1).Value = 45972
$ws.Cells.Item(214,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214,2).Value = 155700
$ws.Cells.Item(214,3).Value = 651102
$ws.Cells.Item(214,4).Value = 469098
$ws.Cells.Item(214,5).Value = 11349509
$ws.Cells.Item(214,6).Value = 73038558600
$ws.Cells.Item(214,7).Value = 1767118551300
$ws.Cells.Item(214,8).Value = 4.130000114440918
$ws.Cells.Item(215,1).Value = 45973
$ws.Cells.Item(215,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215,2).Value = 158000
$ws.Cells.Item(215,3).Value = 213802
$ws.Cells.Item(215,4).Value = 461827
$ws.Cells.Item(215,5).Value = 11349509
$ws.Cells.Item(215,6).Value = 72968666000
$ws.Cells.Item(215,7).Value = 1793222422000
$ws.Cells.Item(215,8).Value = 4.070000171661377
$ws.Cells.Item(216,1).Value = 45974
$ws.Cells.Item(216,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216,2).Value = 160800
$ws.Cells.Item(216,3).Value = 144321
$ws.Cells.Item(217,1).Value = 45975
$ws.Cells.Item(217,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217,2).Value = 156200
$ws.Cells.Item(217,3).Value = 99824

$ws = $wb.Worksheets.Item("에이피알")
$ws.Cells.Item(210,4).Value = 102162
$ws.Cells.Item(210,5).Value = 37430655
$ws.Cells.Item(210,6).Value = 26715363000
$ws.Cells.Item(210,7).Value = 9788116282500
$ws.Cells.Item(210,8).Value = 0.2700000107288361
$ws.Cells.Item(211,4).Value = 139511
$ws.Cells.Item(211,5).Value = 37430655
$ws.Cells.Item(211,6).Value = 32645574000
$ws.Cells.Item(211,7).Value = 8758773270000
$ws.Cells.Item(211,8).Value = 0.3700000047683716
$ws.Cells.Item(212,1).Value = 45968
$ws.Cells.Item(212,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212,2).Value = 234000
$ws.Cells.Item(212,3).Value = 1031039
$ws.Cells.Item(212,4).Value = 145613
$ws.Cells.Item(212,5).Value = 37430655
$ws.Cells.Item(212,6).Value = 34073442000
$ws.Cells.Item(212,7).Value = 8758773270000
$ws.Cells.Item(212,8).Value = 0.3899999856948853
$ws.Cells.Item(213,1).Value = 45971
$ws.Cells.Item(213,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,2).Value = 220500
$ws.Cells.Item(213,3).Value = 919703
$ws.Cells.Item(213,4).Value = 160840
$ws.Cells.Item(213,5).Value = 37430655
$ws.Cells.Item(213,6).Value = 35465220000
$ws.Cells.Item(213,7).Value = 8253459427500
$ws.Cells.Item(213,8).Value = 0.4300000071525574
$ws.Cells.Item(214,1).Value = 45972
$ws.Cells.Item(214,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214,2).Value = 212500
$ws.Cells.Item(214,3).Value = 1269534
$ws.Cells.Item(214,4).Value = 145220
$ws.Cells.Item(214,5).Value = 37430655
$ws.Cells.Item(214,6).Value = 30859250000
$ws.Cells.Item(214,7).Value = 7954014187500
$ws.Cells.Item(214,8).Value = 0.3899999856948853
$ws.Cells.Item(215,1).Value = 45973
$ws.Cells.Item(215,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215,2).Value = 211000
$ws.Cells.Item(215,3).Value = 580602
$ws.Cells.Item(215,4).Value = 112050
$ws.Cells.Item(215,5).Value = 37430655
$ws.Cells.Item(215,6).Value = 23642550000
$ws.Cells.Item(215,7).Value = 7897868205000
$ws.Cells.Item(215,8).Value = 0.300000011920929
$ws.Cells.Item(216,1).Value = 45974
$ws.Cells.Item(216,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216,2).Value = 213500
$ws.Cells.Item(216,3).Value = 473309
$ws.Cells.Item(217,1).Value = 45975
$ws.Cells.Item(217,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217,2).Value = 219500
$ws.Cells.Item(217,3).Value = 703586

$ws = $wb.Worksheets.Item("달바글로벌")
$ws.Cells.Item(113,4).Value = 38472
$ws.Cells.Item(113,5).Value = 12343915
$ws.Cells.Item(113,6).Value = 6190144800
$ws.Cells.Item(113,7).Value = 1986135923500
$ws.Cells.Item(113,8).Value = 0.3100000023841858
$ws.Cells.Item(114,4).Value = 118867
$ws.Cells.Item(114,5).Value = 12343915
$ws.Cells.Item(114,6).Value = 15203089300
$ws.Cells.Item(114,7).Value = 1578786728500
$ws.Cells.Item(114,8).Value = 0.9599999785423279
$ws.Cells.Item(115,1).Value = 45968
$ws.Cells.Item(115,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115,2).Value = 121900
$ws.Cells.Item(115,3).Value = 257902
$ws.Cells.Item(115,4).Value = 108419
$ws.Cells.Item(115,5).Value = 12343915
$ws.Cells.Item(115,6).Value = 13216276100
$ws.Cells.Item(115,7).Value = 1504723238500
$ws.Cells.Item(115,8).Value = 0.8799999952316284
$ws.Cells.Item(116,1).Value = 45971
$ws.Cells.Item(116,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116,2).Value = 121200
$ws.Cells.Item(116,3).Value = 248680
$ws.Cells.Item(116,4).Value = 127352
$ws.Cells.Item(116,5).Value = 12343915
$ws.Cells.Item(116,6).Value = 15435062400
$ws.Cells.Item(116,7).Value = 1496082498000
$ws.Cells.Item(116,8).Value = 1.029999971389771
$ws.Cells.Item(117,1).Value = 45972
$ws.Cells.Item(117,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117,2).Value = 122100
$ws.Cells.Item(117,3).Value = 216968
$ws.Cells.Item(117,4).Value = 136888
$ws.Cells.Item(117,5).Value = 12343915
$ws.Cells.Item(117,6).Value = 16714024800
$ws.Cells.Item(117,7).Value = 1507192021500
$ws.Cells.Item(117,8).Value = 1.110000014305115
$ws.Cells.Item(118,1).Value = 45973
$ws.Cells.Item(118,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(118,2).Value = 126300
$ws.Cells.Item(118,3).Value = 198182
$ws.Cells.Item(118,4).Value = 147169
$ws.Cells.Item(118,5).Value = 12343915
$ws.Cells.Item(118,6).Value = 18587444700
$ws.Cells.Item(118,7).Value = 1559036464500
$ws.Cells.Item(118,8).Value = 1.190000057220459
$ws.Cells.Item(119,1).Value = 45974
$ws.Cells.Item(119,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119,2).Value = 127800
$ws.Cells.Item(119,3).Value = 91769
$ws.Cells.Item(120,1).Value = 45975
$ws.Cells.Item(120,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(120,2).Value = 127900
$ws.Cells.Item(120,3).Value = 106081
